# Milionária draw-history worksheet: append the latest contest results
# (draws 325-328, stored in column A as the 1-based "next row" sequence
# 325..328) to the bottom of the table on "+ MILIONÁRIA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each inner array is one drawing: sequence#, six main numbers, two "trevo" numbers
$newDraws = @(
    @(325, 13, 14, 19, 21, 26, 41, 4, 5),
    @(326, 16, 20, 29, 34, 35, 42, 3, 5),
    @(327,  4,  7, 12, 25, 28, 46, 1, 2),
    @(328,  3, 19, 21, 25, 27, 42, 2, 3)
)

$firstNewRow = 326
for ($i = 0; $i -lt $newDraws.Count; $i++) {
    $rowValues = $newDraws[$i]
    $targetRow = $firstNewRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($targetRow, $col).Value2 = $rowValues[$col - 1]
    }
}

# The previous "last 5 rows" highlight (rows 321:325) moves down to the
# new last 5 rows; drop the leftover formatting from the old ones.
$ws.Range("A321:I325").ClearFormats()

# Leave the cursor where the user would naturally end up after typing
# the new data block, with the view scrolled to keep the table in sight.
$wb.Windows.Item(1).ScrollRow = 297
$ws.Range("D336").Select()
